# add dublin surprise songs
# -----------------------------------------------------------------------
# This script reproduces the "add dublin surprise songs" commit against
# the taylor surprise-songs workbook:
#
#   1. The TTPD-era dress-color labels used in column E (dress) for the
#      already-recorded European shows (rows 168-213, Paris..London) were
#      renamed:
#         pink      -> bright pink
#         orange    -> sunset
#         new blue  -> ocean blue
#   2. The two new Dublin, Ireland show nights (rows 214-219, currently
#      only carrying leg/date/city/night/instrument) get their dress,
#      song and mashup data filled in.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Rename the dress-color shared strings for rows 168-213 --------
$colorRename = @{
    "pink"     = "bright pink"
    "orange"   = "sunset"
    "new blue" = "ocean blue"
}

for ($r = 168; $r -le 213; $r++) {
    $cell = $ws.Cells.Item($r, 5)   # column E = dress
    $cur = $cell.Value2
    if ($colorRename.ContainsKey($cur)) {
        $cell.Value2 = $colorRename[$cur]
    }
}

# ---- 2. Fill in the Dublin, Ireland surprise songs (rows 214-219) -----
# night 1 (2024-09-28), night 2 (2024-09-29), night 3 (2024-09-30)
$ws.Cells.Item(214, 5).Value2 = "ocean blue"
$ws.Cells.Item(214, 7).Value2 = "State Of Grace (Taylor's Version)"
$ws.Cells.Item(214, 8).Value2 = "You're On Your Own, Kid"

$ws.Cells.Item(215, 5).Value2 = "ocean blue"
$ws.Cells.Item(215, 7).Value2 = "Sweet Nothing"
$ws.Cells.Item(215, 8).Value2 = "hoax"

$ws.Cells.Item(216, 5).Value2 = "bright pink"
$ws.Cells.Item(216, 7).Value2 = "The Albatross"
$ws.Cells.Item(216, 8).Value2 = "Dancing With Our Hands Tied"

$ws.Cells.Item(217, 5).Value2 = "bright pink"
$ws.Cells.Item(217, 7).Value2 = "This Love (Taylor's Version)"
$ws.Cells.Item(217, 8).Value2 = "Ours (Taylor's Version)"

$ws.Cells.Item(218, 5).Value2 = "sunset"
$ws.Cells.Item(218, 7).Value2 = "Clara Bow"
$ws.Cells.Item(218, 8).Value2 = "The Lucky One (Taylor's Version)"

$ws.Cells.Item(219, 5).Value2 = "sunset"
$ws.Cells.Item(219, 7).Value2 = "You're On Your Own, Kid"

# ---- 3. Update the saved view state (best-effort, cosmetic only) ------
$ws.Range("F219").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 191
